$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column K into column L for the data rows (4-14),
# then overwrite the values for the new "2023" column.
$ws.Range("K4:K14").Copy($ws.Range("L4:L14"))

$ws.Range("L4").Value = 2023
$ws.Range("L5").Value = 1.6430457248453274
$ws.Range("L6").Value = 0.41181606829870221
$ws.Range("L7").Value = 0.94796963217320562
$ws.Range("L8").Value = 0.72306112208737106
$ws.Range("L9").Value = 2.1802539701246277
$ws.Range("L10").Value = 0.63651150401750112
$ws.Range("L11").Value = 0.97994201681774651
$ws.Range("L12").Value = 2.2469385026996971
$ws.Range("L13").Value = 4.1686356866605365
$ws.Range("L14").Value = 0.3304193846038968

# Row height adjustments
$ws.Rows(2).RowHeight = 13.5
$ws.Rows(3).RowHeight = 13.5

# The saved file no longer pins the cursor on L7; reset the selection to
# the sheet's default top-left cell.
[void]$ws.Range("A1").Select()
